# Update "想去人数" (interest count) values in column F across all four
# sheets, matching the gh-pages data-refresh diff (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1715  # was 1714
$ws.Range("F7").Value = 2763  # was 2762
$ws.Range("F8").Value = 2129  # was 2128
$ws.Range("F9").Value = 881  # was 880
$ws.Range("F10").Value = 2360  # was 2359
$ws.Range("F12").Value = 6847  # was 6846
$ws.Range("F16").Value = 1564  # was 1563
$ws.Range("F17").Value = 1359  # was 1358
$ws.Range("F20").Value = 2811  # was 2806
$ws.Range("F21").Value = 2531  # was 2521
$ws.Range("F22").Value = 2531  # was 2521
$ws.Range("F23").Value = 818  # was 817
$ws.Range("F26").Value = 5481  # was 5479
$ws.Range("F30").Value = 3824  # was 3823
$ws.Range("F33").Value = 1730  # was 1729
$ws.Range("F34").Value = 1091  # was 1090
$ws.Range("F35").Value = 193  # was 191
$ws.Range("F37").Value = 88  # was 87
$ws.Range("F40").Value = 429  # was 428
$ws.Range("F43").Value = 55  # was 52

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 500  # was 499
$ws.Range("F14").Value = 969  # was 968
$ws.Range("F20").Value = 617  # was 616
$ws.Range("F33").Value = 51  # was 50

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F12").Value = 639  # was 637
$ws.Range("F13").Value = 769  # was 766
$ws.Range("F14").Value = 1273  # was 1271

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 1715  # was 1714
$ws.Range("F10").Value = 2763  # was 2762
$ws.Range("F12").Value = 2129  # was 2128
$ws.Range("F13").Value = 881  # was 880
$ws.Range("F14").Value = 2360  # was 2359
$ws.Range("F16").Value = 6847  # was 6846
$ws.Range("F18").Value = 639  # was 637
$ws.Range("F19").Value = 769  # was 766
$ws.Range("F20").Value = 1564  # was 1563
$ws.Range("F21").Value = 1359  # was 1358
$ws.Range("F24").Value = 1273  # was 1271
$ws.Range("F25").Value = 2811  # was 2806
$ws.Range("F26").Value = 2531  # was 2521
$ws.Range("F28").Value = 818  # was 817
$ws.Range("F31").Value = 5481  # was 5479
$ws.Range("F33").Value = 3824  # was 3823
$ws.Range("F36").Value = 1730  # was 1729
$ws.Range("F37").Value = 1091  # was 1090
$ws.Range("F39").Value = 88  # was 87
$ws.Range("F42").Value = 429  # was 428

Write-Output "Applied 47 cell updates"
